$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 55
$ws.Range("I2").Value = 124
$ws.Range("J2").Value = 646
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 169
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 101
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 79
$ws.Range("T2").Value = 102
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 956
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 971
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 1
